$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update competition ID values in column B (rows 2-51) from 66 to 266
$ws.Range("B2:B51").Value = 266
